$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency market snapshot (Price / Volume(1h) columns) with the
# latest scrape, including the Polkadot/WrappedEther (rows 17-18) and
# Stacks/FirstDigitalUSD (rows 43-44) rank swaps. Numeric-looking text values are
# apostrophe-prefixed so Excel keeps them as text (matching this sheet's existing
# text-formatted Price column, e.g. "1.00") instead of coercing them to numbers.

$ws.Range("D2").Value = "58.638.88"
$ws.Range("E2").Value = "  +2.13%  "

$ws.Range("D3").Value = "3.103.51"
$ws.Range("E3").Value = "  +0.65%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'526.76"
$ws.Range("E5").Value = "  +2.25%  "

$ws.Range("D6").Value = "'143.78"
$ws.Range("E6").Value = "  +1.81%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("E8").Value = "  +1.86%  "

$ws.Range("D9").Value = "'7.34"
$ws.Range("E9").Value = "  +1.21%  "

$ws.Range("E10").Value = "  +0.51%  "

$ws.Range("D11").Value = "'0.385"
$ws.Range("E11").Value = "  +3.22%  "

$ws.Range("D12").Value = "3.628.96"
$ws.Range("E12").Value = "  +0.33%  "

$ws.Range("E13").Value = "  +0.82%  "

$ws.Range("D14").Value = "'26.93"
$ws.Range("E14").Value = "  +4.81%  "

$ws.Range("D15").Value = "'0.0000167"
$ws.Range("E15").Value = "  +1.87%  "

$ws.Range("D16").Value = "58.651.35"
$ws.Range("E16").Value = "  +1.90%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "'6.16"
$ws.Range("E17").Value = "  +0.36%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.089.42"
$ws.Range("E18").Value = "  +0.13%  "

$ws.Range("D19").Value = "'12.93"
$ws.Range("E19").Value = "  -1.37%  "

$ws.Range("D20").Value = "'8.12"
$ws.Range("E20").Value = "  -0.19%  "

$ws.Range("D21").Value = "'342.30"
$ws.Range("E21").Value = "  +2.16%  "

$ws.Range("E22").Value = "  -0.11%  "

$ws.Range("D23").Value = "'0.506"
$ws.Range("E23").Value = "  +0.94%  "

$ws.Range("D24").Value = "'65.92"
$ws.Range("E24").Value = "  +0.12%  "

$ws.Range("D25").Value = "'0.171"
$ws.Range("E25").Value = "  +0.17%  "

$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("D27").Value = "0.0₃0918"
$ws.Range("E27").Value = "  -0.74%  "

$ws.Range("D28").Value = "'6.65"
$ws.Range("E28").Value = "  +3.40%  "

$ws.Range("D29").Value = "'7.25"
$ws.Range("E29").Value = "  +1.89%  "

$ws.Range("E30").Value = "  +2.74%  "

$ws.Range("D31").Value = "'21.06"
$ws.Range("E31").Value = "  +1.04%  "

$ws.Range("E32").Value = "  +3.25%  "

$ws.Range("D33").Value = "'154.40"
$ws.Range("E33").Value = "  +0.28%  "

$ws.Range("D34").Value = "'4.66"
$ws.Range("E34").Value = "  +2.37%  "

$ws.Range("D35").Value = "'6.08"
$ws.Range("E35").Value = "  +2.67%  "

$ws.Range("D36").Value = "'27.09"
$ws.Range("E36").Value = "  -2.03%  "

$ws.Range("E37").Value = "  +5.67%  "

$ws.Range("D38").Value = "'0.0678"
$ws.Range("E38").Value = "  +0.03%  "

$ws.Range("D39").Value = "3.140.10"
$ws.Range("E39").Value = "  +0.45%  "

$ws.Range("D40").Value = "'3.91"
$ws.Range("E40").Value = "  +1.55%  "

$ws.Range("D41").Value = "'36.85"
$ws.Range("E41").Value = "  +0.18%  "

$ws.Range("D42").Value = "'0.673"
$ws.Range("E42").Value = "  +0.16%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'1.49"
$ws.Range("E43").Value = "  +6.66%  "

$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  -0.12%  "

$ws.Range("D45").Value = "2.284.29"
$ws.Range("E45").Value = "  -0.31%  "

$ws.Range("E46").Value = "  +1.40%  "

$ws.Range("D47").Value = "'21.03"
$ws.Range("E47").Value = "  +4.49%  "

$ws.Range("D48").Value = "'0.967"
$ws.Range("E48").Value = "  +2.18%  "

$ws.Range("E49").Value = "  +1.80%  "

$ws.Range("D50").Value = "'268.08"
$ws.Range("E50").Value = "  +7.62%  "

$ws.Range("D51").Value = "'0.749"
$ws.Range("E51").Value = "  +8.86%  "
